$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 2024-04-24, now 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the price list values in column D
$ws.Range("D27").Value = 13037.21
$ws.Range("D28").Value = 16298.154
$ws.Range("D29").Value = 13968.909
$ws.Range("D30").Value = 17695.7
$ws.Range("D31").Value = 14900.605
$ws.Range("D32").Value = 19089.968
$ws.Range("D33").Value = 16298.154
$ws.Range("D34").Value = 21894.902
